# Update vm_pu results for Case_4_21 (380 kV case) in row 2-25 (bus 0-23)
# Columns updated per row: B,C,D,E,F,I,J,K,L,M,N (G stays 1, H stays blank)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is: row number, then new values for columns B,C,D,E,F,I,J,K,L,M,N
$updates = @(
    @(2, 1.02, 1.034790442102991, 1.038716206233262, 1.03843578859248, 1.044195948559976, 1.038201528602054, 1.039907758145235, 1.041503447536725, 1.041223828371244, 1.046967679970012, 1.041384546440889),
    @(3, 1.02, 1.035746478966521, 1.03944823678058, 1.039340030694109, 1.045391379218684, 1.038459207182103, 1.04050704043046, 1.042045754640834, 1.041937834376365, 1.047973296641507, 1.041984679775697),
    @(4, 1.02, 1.036365383970425, 1.039922111136398, 1.039925763850363, 1.046165556971385, 1.038624907487465, 1.040894491619613, 1.042396197514024, 1.042399841065588, 1.048624082952428, 1.04237268118998),
    @(5, 1.02, 1.036625639409269, 1.040121375220806, 1.040172155552782, 1.046491178156996, 1.038694319800179, 1.041057297845013, 1.042543411494431, 1.042594067288344, 1.048897693980724, 1.042535718618893),
    @(6, 1.02, 1.036669341402326, 1.040154835283231, 1.040213534544231, 1.046545860592434, 1.03870595987513, 1.041084629121421, 1.042568122771807, 1.042626678624417, 1.048943635679182, 1.04256308870885),
    @(7, 1.02, 1.036368861250651, 1.039924773528183, 1.039929055562072, 1.046169907320934, 1.038625835953429, 1.040896667352553, 1.042398165036094, 1.042402436332197, 1.048627738875843, 1.042374860012711),
    @(8, 1.02, 1.03511347987299, 1.038963557050808, 1.038741250876887, 1.044599814642211, 1.03828882638681, 1.040110355254931, 1.041686818894044, 1.041465130240075, 1.047307515271457, 1.041587431261718),
    @(9, 1.02, 1.032903544229618, 1.037271357737752, 1.036653035782569, 1.041838109732422, 1.037687060066488, 1.038722312146304, 1.040429793753111, 1.039813484744166, 1.044981759725218, 1.040197416972666),
    @(10, 1.02, 1.031431765478655, 1.036144348111608, 1.035264202541089, 1.040000315656323, 1.037280581097304, 1.037795327685763, 1.039589431026478, 1.038712432094599, 1.04343168219719, 1.039269116087864),
    @(11, 1.02, 1.030794832818878, 1.035656618422717, 1.034663617661529, 1.039205318328396, 1.037103317426791, 1.037393554713374, 1.039224995417682, 1.038235682489067, 1.042760579767062, 1.038866772551771),
    @(12, 1.02, 1.030558301682034, 1.035475495785759, 1.034440653007373, 1.038910137518688, 1.037037285417117, 1.037244261317706, 1.039089545248842, 1.038058599054124, 1.042511315823438, 1.038717267142357),
    @(13, 1.02, 1.030609035964145, 1.035514345296639, 1.034488474271253, 1.038973449536755, 1.037051458028356, 1.037276287831052, 1.039118603497303, 1.038096583925349, 1.042564783185549, 1.038749339137025),
    @(14, 1.02, 1.030775279977503, 1.035641645909485, 1.034645184895345, 1.039180916215155, 1.037097863041917, 1.037381215235371, 1.039213800745572, 1.038221044660813, 1.042739975260767, 1.038854415550293),
    @(15, 1.02, 1.030877715568279, 1.03572008559517, 1.034741755315621, 1.039308758717048, 1.037126429746349, 1.037445856929489, 1.039272443988232, 1.038297729388824, 1.042847918643041, 1.038919149043033),
    @(16, 1.02, 1.031474044157355, 1.036176722950189, 1.035304078071485, 1.040053093390813, 1.037292319070591, 1.03782198401286, 1.039613605811267, 1.03874407272828, 1.043476222950506, 1.039295810270003),
    @(17, 1.02, 1.031848201211358, 1.036463233189972, 1.035657020075689, 1.040520203216758, 1.03739604107785, 1.038057816486772, 1.039827459925475, 1.039024056120747, 1.043870365858222, 1.039531977653079),
    @(18, 1.02, 1.032066475137365, 1.036630376038776, 1.035862961259038, 1.040792735734426, 1.037456419206359, 1.03819533663452, 1.039952144020168, 1.039187366956538, 1.044100271586648, 1.039669693095211),
    @(19, 1.02, 1.032140906737907, 1.036687371860388, 1.035933194762522, 1.040885675061722, 1.037476986028895, 1.038242221170872, 1.039994649002161, 1.03924305192747, 1.044178665008594, 1.039716644212981),
    @(20, 1.02, 1.031808054134946, 1.036432490628757, 1.035619144863612, 1.040470078972451, 1.037384925218769, 1.038032517712212, 1.039804520921665, 1.038994016417956, 1.043828077142513, 1.039506642951358),
    @(21, 1.02, 1.030726323765339, 1.035604157888549, 1.034599034222772, 1.039119819235015, 1.037084203120506, 1.037350318300555, 1.039185769803448, 1.038184394000047, 1.042688385194577, 1.038823474738286),
    @(22, 1.02, 1.030046509618546, 1.035083594992471, 1.033958340849843, 1.038271531669202, 1.036894037361513, 1.036921062752169, 1.038796259027296, 1.037675367555903, 1.041971892876106, 1.038393609597786),
    @(23, 1.02, 1.03040686220834, 1.035359531926227, 1.034297918806788, 1.038721161217415, 1.03699495104891, 1.03714865034007, 1.039002791134967, 1.037945210462427, 1.042351711803377, 1.038621520386166),
    @(24, 1.02, 1.031826194764177, 1.036446381787749, 1.03563625880657, 1.0404927277281, 1.037389948371377, 1.038043949253606, 1.039814886235069, 1.039007590062715, 1.043847185565338, 1.039518090726852),
    @(25, 1.02, 1.033474601183874, 1.037708637958999, 1.037192309685262, 1.042551486170152, 1.037843566853255, 1.039081442869872, 1.04075518099247, 1.040240469754542, 1.045582947531306, 1.040557057703055)
)

# Column letters B,C,D,E,F,I,J,K,L,M,N map to these 1-based column indices
$columns = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

foreach ($entry in $updates) {
    $rowNum = $entry[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $colNum = $columns[$i]
        $value = $entry[$i + 1]
        $ws.Cells.Item($rowNum, $colNum).Value = $value
    }
}